$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to remain text while we overwrite values,
# since several values look numeric (e.g. "6.00", "3.73") and a plain
# Value assignment on a General-formatted cell would get coerced to a
# number (losing formatting / becoming a different OOXML cell type).
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '91.389.14'
$ws.Range("E2").Value = '  +0.94%  '
$ws.Range("D3").Value = '3.153.09'
$ws.Range("E3").Value = '  +1.58%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '241.77'
$ws.Range("E5").Value = '  +0.96%  '
$ws.Range("D6").Value = '618.32'
$ws.Range("E6").Value = '  -0.80%  '
$ws.Range("E7").Value = '  +0.91%  '
$ws.Range("D8").Value = '0.375'
$ws.Range("E8").Value = '  +1.41%  '
$ws.Range("E9").Value = '  +0.01%  '
$ws.Range("D10").Value = '3.154.87'
$ws.Range("E10").Value = '  +1.71%  '
$ws.Range("D11").Value = '0.746'
$ws.Range("E11").Value = '  +0.96%  '
$ws.Range("D12").Value = '0.205'
$ws.Range("E12").Value = '  +1.28%  '
$ws.Range("D13").Value = '0.0000249'
$ws.Range("E13").Value = '  -0.14%  '
$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").Value = '35.21'
$ws.Range("E14").Value = '  +0.34%  '
$ws.Range("B15").Value = 'Toncoin'
$ws.Range("C15").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D15").Value = '5.64'
$ws.Range("E15").Value = '  +3.24%  '
$ws.Range("D16").Value = '91.125.02'
$ws.Range("E16").Value = '  +0.76%  '
$ws.Range("D17").Value = '3.740.98'
$ws.Range("D18").Value = '3.120.26'
$ws.Range("E18").Value = '  +1.13%  '
$ws.Range("D19").Value = '3.73'
$ws.Range("E19").Value = '  -2.47%  '
$ws.Range("D20").Value = '15.09'
$ws.Range("E20").Value = '  +5.74%  '
$ws.Range("D21").Value = '6.00'
$ws.Range("E21").Value = '  +5.06%  '
$ws.Range("D22").Value = '458.75'
$ws.Range("E22").Value = '  +2.68%  '
$ws.Range("E23").Value = '  -2.63%  '
$ws.Range("E24").Value = '  +1.70%  '
$ws.Range("D25").Value = '5.93'
$ws.Range("E25").Value = '  +0.57%  '
$ws.Range("D26").Value = '89.08'
$ws.Range("E26").Value = '  -2.36%  '
$ws.Range("D27").Value = '11.88'
$ws.Range("E27").Value = '  -1.09%  '
$ws.Range("B28").Value = 'Hedera'
$ws.Range("C28").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D28").Value = '0.151'
$ws.Range("E28").Value = '  +34.24%  '
$ws.Range("B29").Value = 'WrappedeETH'
$ws.Range("C29").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D29").Value = '3.320.20'
$ws.Range("E29").Value = '  +1.78%  '
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  -0.08%  '
$ws.Range("D31").Value = '0.240'
$ws.Range("E31").Value = '  +10.79%  '
$ws.Range("D32").Value = '0.168'
$ws.Range("E32").Value = '  -6.16%  '
$ws.Range("D33").Value = '9.44'
$ws.Range("E33").Value = '  +2.64%  '
$ws.Range("D34").Value = '0.175'
$ws.Range("E34").Value = '  +12.31%  '
$ws.Range("D35").Value = '26.63'
$ws.Range("E35").Value = '  +0.72%  '
$ws.Range("D36").Value = '7.53'
$ws.Range("E36").Value = '  -0.12%  '
$ws.Range("E37").Value = '  +1.21%  '
$ws.Range("D38").Value = '496.56'
$ws.Range("E38").Value = '  +0.52%  '
$ws.Range("E39").Value = '  +2.65%  '
$ws.Range("D40").Value = '3.87'
$ws.Range("E40").Value = '  -9.98%  '
$ws.Range("D41").Value = '0.449'
$ws.Range("E41").Value = '  +8.11%  '
$ws.Range("D42").Value = '3.44'
$ws.Range("E42").Value = '  -4.38%  '
$ws.Range("E43").Value = '  +0.10%  '
$ws.Range("E44").Value = '  -0.03%  '
$ws.Range("D45").Value = '0.702'
$ws.Range("E45").Value = '  -30.18%  '
$ws.Range("D46").Value = '0.717'
$ws.Range("E46").Value = '  +5.27%  '
$ws.Range("E47").Value = '  +1.52%  '
$ws.Range("D48").Value = '155.22'
$ws.Range("E48").Value = '  -2.93%  '
$ws.Range("D49").Value = '1.37'
$ws.Range("E49").Value = '  +3.01%  '
$ws.Range("D50").Value = '4.54'
$ws.Range("E50").Value = '  +0.21%  '
$ws.Range("D51").Value = '0.0331'
$ws.Range("E51").Value = '  +7.96%  '

# Restore the column to General / default style so the saved XML has no
# residual style index on these cells (matches the original workbook).
$priceRange.NumberFormat = "General"
$priceRange.Style = "Normal"
